$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value for every data row (rows 2-515).
# This automatic update bumps that date by one day (45171 -> 45172) for all rows.
$ws.Range("C2:C515").Value = 45172
